$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Jan"
$ws.Range("A7").Value = "Feb"
$ws.Range("A8").Value = "Mar"
$ws.Range("A9").Value = "Apr"
$ws.Range("A10").Value = "May"
$ws.Range("A11").Value = "Jun"
$ws.Range("A12").Value = "Jul"
$ws.Range("A13").Value = "Aug"
$ws.Range("A14").Value = "Sep"
$ws.Range("A15").Value = "Oct"
$ws.Range("A16").Value = "Nov"
$ws.Range("A17").Value = "Dec"

$ws.Range("A6:A17").Select()
